# Add data for 2022-04-22: updates the "through April 13" week window to
# "through April 14" and bumps the affected neighborhood/month counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new cutoff date.
$ws.Name = "Through 2022-04-14"

# Update the column header text (shared string) to match.
$ws.Range("B1").Value = "April 2022 (through April 14)"

# Updated counts (existing cells changed).
$ws.Range("B2").Value = 7     # Austin / April 2022 (was 3)
$ws.Range("N2").Value = 2     # Austin / April 2019 (was 1)
$ws.Range("R5").Value = 4     # Garfield Park / April 2018 (was 3)
$ws.Range("B8").Value = 3     # Chicago Lawn / April 2022 (was 2)
$ws.Range("F11").Value = 2    # Chatham / April 2021 (was 1)
$ws.Range("B18").Value = 2    # Woodlawn / April 2022 (was 1)

# New counts (previously empty cells).
$ws.Range("V3").Value = 1     # Englewood / April 2017
$ws.Range("N4").Value = 1     # North Lawndale / April 2019
$ws.Range("V6").Value = 1     # Humboldt Park / April 2017
$ws.Range("F14").Value = 1    # Wicker Park / April 2021
$ws.Range("Z18").Value = 1    # Woodlawn / April 2016
$ws.Range("N23").Value = 1    # Auburn Gresham / April 2019
$ws.Range("B28").Value = 1    # West Loop / April 2022
$ws.Range("B40").Value = 1    # Morgan Park / April 2022
$ws.Range("R68").Value = 1    # Fuller Park / April 2018
$ws.Range("B84").Value = 1    # Pullman / April 2022
